# Scheduled-runner update: refresh Marketboard price snapshots (currentAveragePrice*
# and the derived Leve profit columns) across the gathering/crafting Leve sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H5858").Value2 = 1924.5714
$ws.Range("I5858").Value2 = 888.3333
$ws.Range("K5858").Value2 = 2664.9999
$ws.Range("M5858").Value2 = -2514.9999
# Row 62
$ws.Range("H6262").Value2 = 3661.3333
$ws.Range("I6262").Value2 = 3813.3333
$ws.Range("K6262").Value2 = 3813.3333
$ws.Range("M6262").Value2 = -3189.3333
# Row 65
$ws.Range("H6565").Value2 = 3661.3333
$ws.Range("I6565").Value2 = 3813.3333
$ws.Range("K6565").Value2 = 19066.6665
$ws.Range("M6565").Value2 = -15946.6665
# Row 127
$ws.Range("H127127").Value2 = 5499.8335
$ws.Range("I127127").Value2 = 3666.3333
$ws.Range("J127127").Value2 = 7333.3335
$ws.Range("K127127").Value2 = 10998.9999
$ws.Range("L127127").Value2 = 22000.0005
$ws.Range("M127127").Value2 = -6038.999899999999
$ws.Range("N127127").Value2 = -31920.0005
# Row 132
$ws.Range("H132132").Value2 = 76929120
$ws.Range("I132132").Value2 = 100006160
$ws.Range("J132132").Value2 = 5658.3335
$ws.Range("K132132").Value2 = 300018480
$ws.Range("L132132").Value2 = 16975.0005
$ws.Range("M132132").Value2 = -300015950
$ws.Range("N132132").Value2 = -22035.0005
# Row 135
$ws.Range("H135135").Value2 = 1998.3334
$ws.Range("I135135").Value2 = 1498
$ws.Range("J135135").Value2 = 2999
$ws.Range("K135135").Value2 = 13482
$ws.Range("L135135").Value2 = 26991
$ws.Range("M135135").Value2 = -10947
$ws.Range("N135135").Value2 = -32061

# --- ARM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H7474").Value2 = 2488.75
$ws.Range("I7474").Value2 = 2752.5
$ws.Range("J7474").Value2 = 2225
$ws.Range("K7474").Value2 = 2752.5
$ws.Range("L7474").Value2 = 2225
$ws.Range("M7474").Value2 = -1878.5
$ws.Range("N7474").Value2 = -3973
# Row 77
$ws.Range("H7777").Value2 = 2488.75
$ws.Range("I7777").Value2 = 2752.5
$ws.Range("J7777").Value2 = 2225
$ws.Range("K7777").Value2 = 13762.5
$ws.Range("L7777").Value2 = 11125
$ws.Range("M7777").Value2 = -9394.5
$ws.Range("N7777").Value2 = -19861
# Row 132
$ws.Range("H132132").Value2 = 7989.8335
$ws.Range("J132132").Value2 = 3994
$ws.Range("L132132").Value2 = 11982
$ws.Range("N132132").Value2 = -17042

# --- BSM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H8282").Value2 = 20675.75
$ws.Range("I8282").Value2 = 6872.7144
$ws.Range("K8282").Value2 = 6872.7144
$ws.Range("M8282").Value2 = -6489.7144
# Row 85
$ws.Range("H8585").Value2 = 20675.75
$ws.Range("I8585").Value2 = 6872.7144
$ws.Range("K8585").Value2 = 6872.7144
$ws.Range("M8585").Value2 = -5546.7144
# Row 86
$ws.Range("H8686").Value2 = 2808.28
$ws.Range("I8686").Value2 = 2811.8333
$ws.Range("J8686").Value2 = 2799.1428
$ws.Range("K8686").Value2 = 2811.8333
$ws.Range("L8686").Value2 = 2799.1428
$ws.Range("M8686").Value2 = -1688.8333
$ws.Range("N8686").Value2 = -5045.1428
# Row 89
$ws.Range("H8989").Value2 = 2808.28
$ws.Range("I8989").Value2 = 2811.8333
$ws.Range("J8989").Value2 = 2799.1428
$ws.Range("K8989").Value2 = 14059.1665
$ws.Range("L8989").Value2 = 13995.714
$ws.Range("M8989").Value2 = -8443.166499999999
$ws.Range("N8989").Value2 = -25227.714
# Row 94
$ws.Range("H9494").Value2 = 1673.15
$ws.Range("I9494").Value2 = 1768.1177
$ws.Range("J9494").Value2 = 1135
$ws.Range("K9494").Value2 = 1768.1177
$ws.Range("L9494").Value2 = 1135
$ws.Range("M9494").Value2 = -1317.1177
$ws.Range("N9494").Value2 = -2037
# Row 99
$ws.Range("H9999").Value2 = 3608.6365
$ws.Range("J9999").Value2 = 999
$ws.Range("L9999").Value2 = 999
$ws.Range("N9999").Value2 = -3995

# --- CRP sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H4141").Value2 = 15031.6
$ws.Range("I4141").Value2 = 6749.8335
$ws.Range("J4141").Value2 = 20552.777
$ws.Range("K4141").Value2 = 6749.8335
$ws.Range("L4141").Value2 = 20552.777
$ws.Range("M4141").Value2 = -6321.8335
$ws.Range("N4141").Value2 = -21408.777
# Row 50
$ws.Range("H5050").Value2 = 29993.834
$ws.Range("J5050").Value2 = 29993.834
$ws.Range("L5050").Value2 = 29993.834
$ws.Range("N5050").Value2 = -31243.834
# Row 59
$ws.Range("H5959").Value2 = 30388.25
$ws.Range("J5959").Value2 = 34989.08
$ws.Range("L5959").Value2 = 34989.08
$ws.Range("N5959").Value2 = -37279.08
# Row 60
$ws.Range("H6060").Value2 = 22027.072
$ws.Range("I6060").Value2 = 4243
$ws.Range("J6060").Value2 = 24991.084
$ws.Range("K6060").Value2 = 4243
$ws.Range("L6060").Value2 = 24991.084
$ws.Range("M6060").Value2 = -3732
$ws.Range("N6060").Value2 = -26013.084
# Row 68
$ws.Range("H6868").Value2 = 38783.5
$ws.Range("I6868").Value2 = 30268
$ws.Range("K6868").Value2 = 30268
$ws.Range("M6868").Value2 = -29519
# Row 71
$ws.Range("H7171").Value2 = 38783.5
$ws.Range("I7171").Value2 = 30268
$ws.Range("K7171").Value2 = 90804
$ws.Range("M7171").Value2 = -87060
# Row 132
$ws.Range("H132132").Value2 = 1316.1666
$ws.Range("I132132").Value2 = 1316.1666
$ws.Range("K132132").Value2 = 3948.4998
$ws.Range("M132132").Value2 = -1418.4998

# --- CUL sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H44").Value2 = 2143.35
$ws.Range("I44").Value2 = 2158.5
$ws.Range("J44").Value2 = 2133.25
$ws.Range("K44").Value2 = 6475.5
$ws.Range("L44").Value2 = 6399.75
$ws.Range("M44").Value2 = -6363.5
$ws.Range("N44").Value2 = -6623.75
# Row 33
$ws.Range("H3333").Value2 = 76.833336
$ws.Range("I3333").Value2 = 17.75
$ws.Range("K3333").Value2 = 106.5
$ws.Range("M3333").Value2 = 176.5
# Row 80
$ws.Range("H8080").Value2 = 1498
$ws.Range("I8080").Value2 = 1498
$ws.Range("J8080").Value2 = 0
$ws.Range("K8080").Value2 = 4494
$ws.Range("L8080").Value2 = 0
$ws.Range("M8080").Value2 = -3558
$ws.Range("N8080").ClearContents()
# Row 83
$ws.Range("H8383").Value2 = 1498
$ws.Range("I8383").Value2 = 1498
$ws.Range("J8383").Value2 = 0
$ws.Range("K8383").Value2 = 13482
$ws.Range("L8383").Value2 = 0
$ws.Range("M8383").Value2 = -8802
$ws.Range("N8383").ClearContents()
# Row 122
$ws.Range("H122122").Value2 = 1262.8
$ws.Range("I122122").Value2 = 1202.6
$ws.Range("J122122").Value2 = 1323
$ws.Range("K122122").Value2 = 10823.4
$ws.Range("L122122").Value2 = 11907
$ws.Range("M122122").Value2 = -8373.4
$ws.Range("N122122").Value2 = -16807
# Row 131
$ws.Range("H131131").Value2 = 2449.4285
$ws.Range("J131131").Value2 = 3112.25
$ws.Range("L131131").Value2 = 9336.75
$ws.Range("N131131").Value2 = -19416.75

# --- GSM sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113113").Value2 = 2057.3
$ws.Range("I113113").Value2 = 1851
$ws.Range("K113113").Value2 = 1851
$ws.Range("M113113").Value2 = 319

# --- LTW sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H9393").Value2 = 642.8570999999999
$ws.Range("I9393").Value2 = 642.8570999999999
$ws.Range("K9393").Value2 = 642.8570999999999
$ws.Range("M9393").Value2 = 605.1429000000001
# Row 122
$ws.Range("H122122").Value2 = 3001.111
$ws.Range("I122122").Value2 = 3001.111
$ws.Range("K122122").Value2 = 9003.332999999999
$ws.Range("M122122").Value2 = -6553.332999999999

# --- WVR sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H8181").Value2 = 533
$ws.Range("I8181").Value2 = 533
$ws.Range("K8181").Value2 = 1066
$ws.Range("M8181").Value2 = -5
# Row 84
$ws.Range("H8484").Value2 = 533
$ws.Range("I8484").Value2 = 533
$ws.Range("K8484").Value2 = 5330
$ws.Range("M8484").Value2 = -26
